$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these text cells remain stored as text (not auto-converted to numbers/percentages)
# by explicitly setting the NumberFormat to Text ("@") before assigning the new values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.02%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.89"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.87%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.914"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.22%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08302"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.65%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.788"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.18%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.499"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.01%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.954"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.42%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.25%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9309"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.57%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1253"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.16%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1940"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.14%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09490"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.73%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03971"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "6.71%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1065"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.12%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001300"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.04%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006075"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.56%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.529"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.77%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.155"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "9.77%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1372"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.63%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2600"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.00%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04423"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.59%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001258"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.93%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004410"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.51%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001193"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.90%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.16%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02831"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.89%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05637"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.82%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007933"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.52%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.55%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009089"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.95%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002106"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.42%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008797"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-28.67%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007308"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "6.27%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003655"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "3.53%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.02%"
